$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" (strikeouts) values replacing the old "Strike#" values in column G,
# keyed by worksheet row number.
$kValues = @{
    2  = 1
    3  = 2
    4  = 2
    5  = 2
    6  = 0
    7  = 1
    8  = 0
    9  = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 1
    14 = 1
    15 = 2
    17 = 0
    18 = 0
    19 = 0
    20 = 2
    21 = 1
    22 = 2
    23 = 1
    24 = 0
    25 = 0
    27 = 1
    28 = 0
    29 = 0
    30 = 3
    31 = 2
    32 = 0
    33 = 0
    34 = 1
    35 = 0
    36 = 1
    37 = 0
    38 = 0
    39 = 0
    40 = 1
    41 = 1
    42 = 1
    43 = 1
    44 = 1
    45 = 1
    46 = 1
    47 = 1
    48 = 1
    49 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
